# "reverted back to summing word fmris to get results"
#
# 1. The three most-recent slurm job notes (P69:P71) are updated to the new
#    job ids that were produced after reverting to the "summed word-level
#    fmris" approach.
# 2. A cautionary note is added (column Q) to the six rows that used the
#    no-overlap word-split evaluation, explaining that it somewhat inflated
#    results.
# 3. A new results row (72) is appended, recording a "local" run of the
#    reverted (summed word fmris) approach.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Update the slurm job ids for the three most recent runs ---------
$ws.Cells.Item(69, 16).Value = "slurm-42642165"
$ws.Cells.Item(70, 16).Value = "slurm-42642207"
$ws.Cells.Item(71, 16).Value = "slurm-42642224"

# --- 2. Add the caveat note to rows 57-62 (column Q) ---------------------
$note = "note: these assements had the text split by the number of words so there was no overlap, thus somewhat inflating results"
$ws.Cells.Item(57, 17).Value = $note
$ws.Cells.Item(58, 17).Value = $note
$ws.Cells.Item(59, 17).Value = $note
$ws.Cells.Item(60, 17).Value = $note
$ws.Cells.Item(61, 17).Value = $note
$ws.Cells.Item(62, 17).Value = $note

# --- 3. Append new row 72 -------------------------------------------------
$ws.Cells.Item(72, 1).Value = "leave-one-out average"
$ws.Cells.Item(72, 2).Value = "4 words"
$ws.Cells.Item(72, 3).Value = "(fmri channel for each word) (detrended) gaussian weighted 3D fMRI image 2-8 seconds after each word"
$ws.Cells.Item(72, 4).Value = "Default CLIP"
$ws.Cells.Item(72, 5).Value = "Default CLIP"
$ws.Cells.Item(72, 6).Value = "3D Resnet18"
$ws.Cells.Item(72, 7).Value = "Cosine Similarity"
$ws.Cells.Item(72, 8).Value = "embed_dim=1024, image_resolution, layers=(2,2,2,2), width=64, context_length=16, vocab_size, transformer_width, transformer_heads, transformer_layers"
$ws.Cells.Item(72, 9).Value = "LR=1e-5, batch_size=32, weight_decay=0.2"
$ws.Cells.Item(72, 10).Value = 50
$ws.Cells.Item(72, 11).Value = "700/0/100"
$ws.Cells.Item(72, 16).Value = "local"

# --- Update the active selection to mirror the author's cursor position --
$ws.Activate() | Out-Null
$ws.Range("P72").Select() | Out-Null
